# Generate Report for Handoff
# Replace the old handoff id (f094f450-ceda-45c9-982c-79d073efdd0c) with the
# new one (db348876-a191-4563-82dc-a4e16d5055d1) everywhere it appears, and
# bump the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldId = "f094f450-ceda-45c9-982c-79d073efdd0c"
$newId = "db348876-a191-4563-82dc-a4e16d5055d1"

$newMd  = "$newId.md"

$newRelMd = "e2e\$newId.md"

$newZhXlf = "$newId.367f1bc4af17f40373469def0393d72f57474047.zh-cn.xlf"

$newDeXlf = "$newId.367f1bc4af17f40373469def0393d72f57474047.de-de.xlf"

$newHoDate   = "2016-12-15 04:38:56"
$newZhXlfDate = "2016-12-15 04:38:43"

# The hyperlink relationship targets (xl/worksheets/_rels/*.rels) are left
# untouched by this change, so the underlying hyperlink Address must stay
# pointed at the original (old-id) URL - only the visible display text
# changes.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fcf5bf5e23d41dc7cbf53ae093fbad0de4d2c6a6/e2e/$oldId.md"

# ---------------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Value = $newRelMd
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkUrl, "", "", $newRelMd)

# ---------------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMd
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhXlfDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkUrl, "", "", $newMd)

# ---------------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMd
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHoDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkUrl, "", "", $newMd)
